# Auto-generated COM-interop edit script
# Updates market-derived profit columns (H:N) for specific leve rows
# across several job sheets, matching a scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 69.583336
$ws.Range("I9").Value = 60.625
$ws.Range("J9").Value = 87.5
$ws.Range("K9").Value = 60.625
$ws.Range("L9").Value = 87.5
$ws.Range("M9").Value = 108.375
$ws.Range("N9").Value = -425.5

# Row 28
$ws.Range("H28").Value = 1313.125
$ws.Range("I28").Value = 1067.3334
$ws.Range("J28").Value = 5000
$ws.Range("K28").Value = 1067.3334
$ws.Range("L28").Value = 5000
$ws.Range("M28").Value = -582.3334
$ws.Range("N28").Value = -5970

# Row 33
$ws.Range("H33").Value = 1079313.2
$ws.Range("J33").Value = 4066.6667
$ws.Range("L33").Value = 4066.6667
$ws.Range("N33").Value = -4524.6667

# Row 51
$ws.Range("H51").Value = 5558662.5
$ws.Range("I51").Value = 3371.5715
$ws.Range("J51").Value = 10419542
$ws.Range("K51").Value = 3371.5715
$ws.Range("L51").Value = 10419542
$ws.Range("M51").Value = -2887.5715
$ws.Range("N51").Value = -10420510

# Row 96
$ws.Range("H96").Value = 260423.5
$ws.Range("I96").Value = 402360.6
$ws.Range("J96").Value = 23861.666
$ws.Range("K96").Value = 1207081.8
$ws.Range("L96").Value = 71584.99800000001
$ws.Range("M96").Value = -1205708.8
$ws.Range("N96").Value = -74330.99800000001

# Row 106
$ws.Range("H106").Value = 57931.6
$ws.Range("I106").Value = 79462.16
$ws.Range("K106").Value = 79462.16
$ws.Range("M106").Value = -78831.16

# Row 116
$ws.Range("H116").Value = 40029.09
$ws.Range("I116").Value = 16333
$ws.Range("J116").Value = 48915.125
$ws.Range("K116").Value = 16333
$ws.Range("L116").Value = 48915.125
$ws.Range("M116").Value = -12891
$ws.Range("N116").Value = -55799.125

# Row 132
$ws.Range("H132").Value = 3445.5745
$ws.Range("I132").Value = 2445.5366
$ws.Range("K132").Value = 7336.6098
$ws.Range("M132").Value = -4806.6098

# Row 138
$ws.Range("H138").Value = 3536.1372
$ws.Range("I138").Value = 1795.9375
$ws.Range("J138").Value = 4331.657
$ws.Range("K138").Value = 5387.8125
$ws.Range("L138").Value = 12994.971
$ws.Range("M138").Value = -247.8125
$ws.Range("N138").Value = -23274.971

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1602.6154
$ws.Range("I2").Value = 1509.3158
$ws.Range("K2").Value = 1509.3158
$ws.Range("M2").Value = -1396.3158

# Row 32
$ws.Range("H32").Value = 12286559
$ws.Range("I32").Value = 6557455.5
$ws.Range("J32").Value = 29473872
$ws.Range("K32").Value = 6557455.5
$ws.Range("L32").Value = 29473872
$ws.Range("M32").Value = -6557168.5
$ws.Range("N32").Value = -29474446

# Row 74
$ws.Range("H74").Value = 1975.8611
$ws.Range("I74").Value = 1638.4814
$ws.Range("J74").Value = 2988
$ws.Range("K74").Value = 1638.4814
$ws.Range("L74").Value = 2988
$ws.Range("M74").Value = -764.4813999999999
$ws.Range("N74").Value = -4736

# Row 77
$ws.Range("H77").Value = 1975.8611
$ws.Range("I77").Value = 1638.4814
$ws.Range("J77").Value = 2988
$ws.Range("K77").Value = 8192.406999999999
$ws.Range("L77").Value = 14940
$ws.Range("M77").Value = -3824.406999999999
$ws.Range("N77").Value = -23676

# Row 97
$ws.Range("H97").Value = 37545.043
$ws.Range("J97").Value = 87917.8
$ws.Range("L97").Value = 87917.8
$ws.Range("N97").Value = -88909.8

# Row 113
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

# Row 116
$ws.Range("H116").Value = 1602.6154
$ws.Range("I116").Value = 1509.3158
$ws.Range("K116").Value = 1509.3158
$ws.Range("M116").Value = 784.6841999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1602.6154
$ws.Range("I3").Value = 1509.3158
$ws.Range("K3").Value = 1509.3158
$ws.Range("M3").Value = -1395.3158

# Row 86
$ws.Range("H86").Value = 1580.4
$ws.Range("J86").Value = 1732.3077
$ws.Range("L86").Value = 1732.3077
$ws.Range("N86").Value = -3978.3077

# Row 88
$ws.Range("H88").Value = 41244.25
$ws.Range("I88").Value = 7500
$ws.Range("J88").Value = 44311.91
$ws.Range("K88").Value = 7500
$ws.Range("L88").Value = 44311.91
$ws.Range("N88").Value = -45123.91
$ws.Range("M88").Value = -7094

# Row 89
$ws.Range("H89").Value = 1580.4
$ws.Range("J89").Value = 1732.3077
$ws.Range("L89").Value = 8661.538500000001
$ws.Range("N89").Value = -19893.5385

# Row 91
$ws.Range("H91").Value = 41244.25
$ws.Range("I91").Value = 7500
$ws.Range("J91").Value = 44311.91
$ws.Range("K91").Value = 7500
$ws.Range("L91").Value = 44311.91
$ws.Range("N91").Value = -47119.91
$ws.Range("M91").Value = -6096

# Row 99
$ws.Range("H99").Value = 47667.734
$ws.Range("J99").Value = 169366.5
$ws.Range("L99").Value = 169366.5
$ws.Range("N99").Value = -172362.5

$ws = $wb.Worksheets.Item("CRP")
# Row 134
$ws.Range("H134").Value = 2110.48
$ws.Range("I134").Value = 1941.6171
$ws.Range("J134").Value = 4756
$ws.Range("K134").Value = 5824.8513
$ws.Range("L134").Value = 14268
$ws.Range("M134").Value = -3289.8513
$ws.Range("N134").Value = -19338

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 70.25
$ws.Range("I2").Value = 22.2
$ws.Range("K2").Value = 133.2
$ws.Range("M2").Value = -20.19999999999999

# Row 68
$ws.Range("H68").Value = 864163.1
$ws.Range("J68").Value = 1390869.6
$ws.Range("L68").Value = 4172608.8
$ws.Range("N68").Value = -4174230.8

# Row 71
$ws.Range("H71").Value = 864163.1
$ws.Range("J71").Value = 1390869.6
$ws.Range("L71").Value = 12517826.4
$ws.Range("N71").Value = -12525938.4

# Row 107
$ws.Range("H107").Value = 2858.125
$ws.Range("J107").Value = 3123.5715
$ws.Range("L107").Value = 9370.7145
$ws.Range("N107").Value = -13210.7145

# Row 137
$ws.Range("H137").Value = 4548795
$ws.Range("I137").Value = 9092170
$ws.Range("J137").Value = 5419.909
$ws.Range("K137").Value = 27276510
$ws.Range("L137").Value = 16259.727
$ws.Range("M137").Value = -27271410
$ws.Range("N137").Value = -26459.727

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 43293796
$ws.Range("I70").Value = 5069.2
$ws.Range("J70").Value = 86582520
$ws.Range("K70").Value = 5069.2
$ws.Range("L70").Value = 86582520
$ws.Range("M70").Value = -4799.2
$ws.Range("N70").Value = -86583060

# Row 73
$ws.Range("H73").Value = 43293796
$ws.Range("I73").Value = 5069.2
$ws.Range("J73").Value = 86582520
$ws.Range("K73").Value = 5069.2
$ws.Range("L73").Value = 86582520
$ws.Range("M73").Value = -4133.2
$ws.Range("N73").Value = -86584392

$ws = $wb.Worksheets.Item("WVR")
# Row 93
$ws.Range("H93").Value = 28213.4
$ws.Range("J93").Value = 28213.4
$ws.Range("L93").Value = 28213.4
$ws.Range("N93").Value = -33205.4

